$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1/G1, copy formatting from E1 (header style)
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "Type of type"
$ws.Range("G1").Value = "Level of type"

# New row 5 data
# A5/D5 date-style cells -> copy format from A4/D4
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# E5 text-style cell -> copy format from E4
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats

# C5/F5/G5 default-style cells -> copy format from C2 (default style)
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 44001.8541666667
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 44001.84375
$ws.Range("E5").Value = "Playing Fortnite"
$ws.Range("F5").Value = "P"
$ws.Range("G5").Value = 1

# Selection moves to A6
$ws.Range("A6").Select() | Out-Null

$excel.CutCopyMode = $false
